# Apply the workbook edit described by the commit:
#  1. Rename two shared-string headers: ht_goals_h -> HTHG, ht_goals_a -> HTAG
#  2. Swap the data (columns B..AD) between specific pairs of rows
#     (the source rows were interleaved / reordered in the refreshed
#     scrape, while the running index in column A stayed put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header rename -------------------------------------------------
$ws.Range("I1").Value = "HTHG"
$ws.Range("J1").Value = "HTAG"

# --- 2. Row-pair content swaps (columns B through AD) -----------------
$rowPairs = @(
    @(30, 31),
    @(66, 67),
    @(70, 71),
    @(94, 95),
    @(97, 98),
    @(102, 103),
    @(128, 129),
    @(133, 134),
    @(151, 152),
    @(164, 165),
    @(177, 178),
    @(238, 239),
    @(267, 268),
    @(287, 288),
    @(291, 292),
    @(296, 297),
    @(301, 303),
    @(302, 304),
    @(305, 306)
)

$firstCol = 2   # column B
$lastCol  = 30  # column AD

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $c1 = $ws.Cells.Item($r1, $col)
        $c2 = $ws.Cells.Item($r2, $col)

        $v1 = $c1.Value2
        $v2 = $c2.Value2

        $c1.Value = $v2
        $c2.Value = $v1
    }
}
